$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B74: change from text "2" to numeric 2
$ws.Cells.Item(74, 2).Value = 2

# Add new row 75
$ws.Cells.Item(75, 1).Value = "Ruilin"

# B75 keeps politeness_score stored as text "3" (not numeric), matching source data
$ws.Cells.Item(75, 2).NumberFormat = "@"
$ws.Cells.Item(75, 2).Value = "3"
$ws.Cells.Item(75, 2).ClearFormats()

$ws.Cells.Item(75, 3).Value = "无"
$ws.Cells.Item(75, 4).Value = "DFT"
$ws.Cells.Item(75, 5).Value = "MET"
$ws.Cells.Item(75, 6).Value = "bdd42601-fca4-49a7-8203-fa53b228c875"
$ws.Cells.Item(75, 7).Value = "BkA7gfZAb_annotated.xlsx"
$ws.Cells.Item(75, 8).Value = "Note that d(A, B'_theta) is *equal* to min_alpha max_w (...) above equation (2) (it is not just an upper bound)."
